$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert two new columns before column D (shifts D:K -> F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Step 2: Copy number formats from the shifted-original D:E range (now F:G) into new D:E
# so the new quarter columns inherit the same date/number styles as the rest of the table.
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Step 2b: A few rows have no data in columns D:M at all (section headers / separators);
# remove the blank formatted cells PasteSpecial introduced there so the row stays empty.
$ws.Range("D5:E6").Clear()
$ws.Range("D36:E37").Clear()
$ws.Range("D78:E79").Clear()

# Step 3: Populate the two new quarter columns (D = 2018-12-31, E = 2018-09-30) with data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 14411000
$ws.Range("E8").Value = 14318000
$ws.Range("D9").Value = 12469000
$ws.Range("E9").Value = 12397000
$ws.Range("D10").Value = 1942000
$ws.Range("E10").Value = 1921000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 12560000
$ws.Range("E17").Value = 12355000
$ws.Range("D18").Value = 1851000
$ws.Range("E18").Value = 1963000
$ws.Range("D20").Value = -197000
$ws.Range("E20").Value = -211000
$ws.Range("D21").Value = 1958000
$ws.Range("E21").Value = 2043000
$ws.Range("D22").Value = 171000
$ws.Range("E22").Value = 177000
$ws.Range("D23").Value = 1483000
$ws.Range("E23").Value = 1575000
$ws.Range("D24").Value = 273000
$ws.Range("E24").Value = 102000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 1210000
$ws.Range("E26").Value = 1473000
$ws.Range("D27").Value = 1210000
$ws.Range("E27").Value = 1473000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 43000
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 197000
$ws.Range("E32").Value = 211000
$ws.Range("D33").Value = 1253000
$ws.Range("E33").Value = 1473000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 1253000
$ws.Range("E35").Value = 1473000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 772000
$ws.Range("E41").Value = 897000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 11916000
$ws.Range("E43").Value = 12185000
$ws.Range("D44").Value = 2997000
$ws.Range("E44").Value = 3050000
$ws.Range("D45").Value = 418000
$ws.Range("E45").Value = 727000
$ws.Range("D46").Value = 16103000
$ws.Range("E46").Value = 16859000
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 6124000
$ws.Range("E48").Value = 5902000
$ws.Range("D49").Value = 14263000
$ws.Range("E49").Value = 14358000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 8386000
$ws.Range("E52").Value = 8376000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 44876000
$ws.Range("E54").Value = 45495000
$ws.Range("D57").Value = 2402000
$ws.Range("E57").Value = 2691000
$ws.Range("D58").Value = 1500000
$ws.Range("E58").Value = 1240000
$ws.Range("D59").Value = 10496000
$ws.Range("E59").Value = 11273000
$ws.Range("D60").Value = 14398000
$ws.Range("E60").Value = 15204000
$ws.Range("D61").Value = 12604000
$ws.Range("E61").Value = 13486000
$ws.Range("D62").Value = 16425000
$ws.Range("E62").Value = 15803000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 43482000
$ws.Range("E66").Value = 44552000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 15434000
$ws.Range("E72").Value = 14737000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1394000
$ws.Range("E76").Value = 943000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 1253000
$ws.Range("E81").Value = 1473000
$ws.Range("D83").Value = 304000
$ws.Range("E83").Value = 291000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 2217000
$ws.Range("E89").Value = 361000
$ws.Range("D91").Value = -459000
$ws.Range("E91").Value = -339000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -402000
$ws.Range("E94").Value = -344000
$ws.Range("D96").Value = -622000
$ws.Range("E96").Value = -569000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -1940000
$ws.Range("E100").Value = -301000
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -125000
$ws.Range("E102").Value = -284000

# Step 4: Apply restated figures for the quarters that moved into columns H/I (and I/J for row 91)
$ws.Range("H8").Value = 13844000
$ws.Range("I8").Value = 12341000
$ws.Range("H9").Value = 12135000
$ws.Range("I9").Value = 10741000
$ws.Range("H10").Value = 1709000
$ws.Range("I10").Value = 1600000
$ws.Range("H17").Value = 11895000
$ws.Range("I17").Value = 10664000
$ws.Range("H18").Value = 1949000
$ws.Range("I18").Value = 1677000
$ws.Range("H20").Value = -203000
$ws.Range("I20").Value = -218000
$ws.Range("H21").Value = 2061000
$ws.Range("I21").Value = 1758000
$ws.Range("H23").Value = 1572000
$ws.Range("I23").Value = 1297000
$ws.Range("H24").Value = 503000
$ws.Range("I24").Value = 334000
$ws.Range("H26").Value = 1069000
$ws.Range("I26").Value = 963000
$ws.Range("H27").Value = 1069000
$ws.Range("I27").Value = 963000
$ws.Range("H29").Value = -1813000
$ws.Range("H32").Value = 203000
$ws.Range("I32").Value = 218000
$ws.Range("H33").Value = -744000
$ws.Range("I33").Value = 963000
$ws.Range("H35").Value = -744000
$ws.Range("I35").Value = 963000
$ws.Range("H81").Value = -744000
$ws.Range("I81").Value = 963000
$ws.Range("I91").Value = -222000
$ws.Range("J91").Value = -278000
